$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.685.35"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "3.474.27"
$ws.Range("E3").Value = "  +4.74%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'262.50"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").Value = "'675.87"
$ws.Range("E6").Value = "  +8.39%  "
$ws.Range("E7").Value = "  +9.17%  "
$ws.Range("D8").Value = "'0.462"
$ws.Range("E8").Value = "  +14.35%  "
$ws.Range("D9").Value = "'1.13"
$ws.Range("E9").Value = "  +24.13%  "
$ws.Range("D10").Value = "'0.998"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.472.61"
$ws.Range("E11").Value = "  +4.68%  "
$ws.Range("E12").Value = "  +9.93%  "
$ws.Range("D13").Value = "'43.50"
$ws.Range("E13").Value = "  +11.07%  "
$ws.Range("E14").Value = "  +9.93%  "
$ws.Range("D15").Value = "'6.31"
$ws.Range("E15").Value = "  +14.75%  "
$ws.Range("D16").Value = "98.310.37"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "4.115.61"
$ws.Range("E17").Value = "  +4.50%  "
$ws.Range("D18").Value = "'8.82"
$ws.Range("E18").Value = "  +40.22%  "
$ws.Range("D19").Value = "3.480.54"
$ws.Range("D20").Value = "'17.83"
$ws.Range("E20").Value = "  +16.59%  "
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("D22").Value = "'533.35"
$ws.Range("E22").Value = "  +10.26%  "
$ws.Range("D23").Value = "'11.07"
$ws.Range("E23").Value = "  +17.31%  "
$ws.Range("D24").Value = "'0.492"
$ws.Range("E24").Value = "  +67.87%  "
$ws.Range("E25").Value = "  +6.50%  "
$ws.Range("D26").Value = "'6.51"
$ws.Range("E26").Value = "  +15.97%  "
$ws.Range("D27").Value = "'103.98"
$ws.Range("E27").Value = "  +17.20%  "
$ws.Range("D28").Value = "'13.18"
$ws.Range("E28").Value = "  +9.98%  "
$ws.Range("D29").Value = "'0.154"
$ws.Range("E29").Value = "  +17.26%  "
$ws.Range("D30").Value = "'11.84"
$ws.Range("E30").Value = "  +15.75%  "
$ws.Range("D31").Value = "'0.198"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'0.602"
$ws.Range("E33").Value = "  +29.92%  "
$ws.Range("D34").Value = "'30.66"
$ws.Range("E34").Value = "  +9.54%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +16.28%  "
$ws.Range("D37").Value = "'8.18"
$ws.Range("E37").Value = "  +12.62%  "
$ws.Range("E38").Value = "  +8.75%  "
$ws.Range("D39").Value = "'539.63"
$ws.Range("E39").Value = "  +9.57%  "
$ws.Range("E40").Value = "  +15.95%  "
$ws.Range("D41").Value = "'24.78"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'0.0446"
$ws.Range("E42").Value = "  +35.78%  "
$ws.Range("E43").Value = "  +9.63%  "
$ws.Range("E44").Value = "  +11.72%  "
$ws.Range("D45").Value = "'3.74"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").Value = "'8.72"
$ws.Range("E46").Value = "  +18.00%  "
$ws.Range("D47").Value = "'1.63"
$ws.Range("E47").Value = "  +20.13%  "
$ws.Range("D48").Value = "'5.43"
$ws.Range("E48").Value = "  +16.74%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +10.89%  "
$ws.Range("D51").Value = "'52.50"
$ws.Range("E51").Value = "  +15.99%  "
